$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking cells so they remain strings like the source data
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(51,4).NumberFormat = "@"

$ws.Cells.Item(2,4).Value = '26.697.23'
$ws.Cells.Item(2,5).Value = '  +1.77%  '

$ws.Cells.Item(3,4).Value = '1.636.83'
$ws.Cells.Item(3,5).Value = '  +2.00%  '

$ws.Cells.Item(4,5).Value = '  -0.06%  '

$ws.Cells.Item(5,4).Value = '213.10'
$ws.Cells.Item(5,5).Value = '  +0.18%  '

$ws.Cells.Item(6,5).Value = '  +2.17%  '

$ws.Cells.Item(7,5).Value = '  -0.07%  '

$ws.Cells.Item(8,5).Value = '  +1.60%  '

$ws.Cells.Item(9,4).Value = '0.0624'
$ws.Cells.Item(9,5).Value = '  +1.82%  '

$ws.Cells.Item(10,4).Value = '19.06'
$ws.Cells.Item(10,5).Value = '  +3.47%  '

$ws.Cells.Item(11,5).Value = '  +2.93%  '

$ws.Cells.Item(12,5).Value = '  +1.91%  '

$ws.Cells.Item(13,4).Value = '1.629.99'
$ws.Cells.Item(13,5).Value = '  +1.56%  '

$ws.Cells.Item(14,4).Value = '4.08'
$ws.Cells.Item(14,5).Value = '  +1.53%  '

$ws.Cells.Item(15,4).Value = '0.527'
$ws.Cells.Item(15,5).Value = '  +2.85%  '

$ws.Cells.Item(16,4).Value = '26.697.06'
$ws.Cells.Item(16,5).Value = '  +1.88%  '

$ws.Cells.Item(17,4).Value = '63.12'
$ws.Cells.Item(17,5).Value = '  +1.75%  '

$ws.Cells.Item(18,5).Value = '  +1.83%  '

$ws.Cells.Item(19,4).Value = '208.97'
$ws.Cells.Item(19,5).Value = '  +3.98%  '

$ws.Cells.Item(20,5).Value = '  +0.02%  '

$ws.Cells.Item(21,4).Value = '4.31'
$ws.Cells.Item(21,5).Value = '  +1.17%  '

$ws.Cells.Item(22,4).Value = '9.40'
$ws.Cells.Item(22,5).Value = '  +1.12%  '

$ws.Cells.Item(23,4).Value = '6.17'
$ws.Cells.Item(23,5).Value = '  +2.83%  '

$ws.Cells.Item(24,5).Value = '  +1.65%  '

$ws.Cells.Item(25,4).Value = '146.35'
$ws.Cells.Item(25,5).Value = '  +1.61%  '

$ws.Cells.Item(26,5).Value = '  -0.11%  '

$ws.Cells.Item(27,5).Value = '  -0.70%  '

$ws.Cells.Item(28,4).Value = '6.74'
$ws.Cells.Item(28,5).Value = '  +2.72%  '

$ws.Cells.Item(29,4).Value = '15.40'
$ws.Cells.Item(29,5).Value = '  +1.28%  '

$ws.Cells.Item(30,4).Value = '0.0521'
$ws.Cells.Item(30,5).Value = '  +6.18%  '

$ws.Cells.Item(31,5).Value = '  -0.15%  '

$ws.Cells.Item(32,4).Value = '3.24'
$ws.Cells.Item(32,5).Value = '  +1.14%  '

$ws.Cells.Item(33,5).Value = '  +0.90%  '

$ws.Cells.Item(34,2).Value = 'LidoDAOToken'
$ws.Cells.Item(34,3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(34,4).Value = '1.50'
$ws.Cells.Item(34,5).Value = '  +1.02%  '

$ws.Cells.Item(35,2).Value = 'HuobiToken'
$ws.Cells.Item(35,3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35,4).Value = '2.41'
$ws.Cells.Item(35,5).Value = '  +1.49%  '

$ws.Cells.Item(36,4).Value = '1.169.81'
$ws.Cells.Item(36,5).Value = '  +0.74%  '

$ws.Cells.Item(37,5).Value = '  -0.89%  '

$ws.Cells.Item(38,4).Value = '0.809'
$ws.Cells.Item(38,5).Value = '  +3.10%  '

$ws.Cells.Item(39,5).Value = '  -0.02%  '

$ws.Cells.Item(40,4).Value = '0.505'
$ws.Cells.Item(40,5).Value = '  +1.69%  '

$ws.Cells.Item(41,5).Value = '  +0.29%  '

$ws.Cells.Item(42,4).Value = '0.794'
$ws.Cells.Item(42,5).Value = '  +1.28%  '

$ws.Cells.Item(43,5).Value = '  +1.10%  '

$ws.Cells.Item(44,4).Value = '1.775.41'
$ws.Cells.Item(44,5).Value = '  +2.08%  '

$ws.Cells.Item(45,4).Value = '92.35'
$ws.Cells.Item(45,5).Value = '  +0.71%  '

$ws.Cells.Item(46,4).Value = '1.56'
$ws.Cells.Item(46,5).Value = '  +1.55%  '

$ws.Cells.Item(47,5).Value = '  -2.43%  '

$ws.Cells.Item(48,4).Value = '54.79'
$ws.Cells.Item(48,5).Value = '  +1.21%  '

$ws.Cells.Item(49,5).Value = '  +1.66%  '

$ws.Cells.Item(50,4).Value = '0.410'
$ws.Cells.Item(50,5).Value = '  +0.75%  '

$ws.Cells.Item(51,4).Value = '7.54'
$ws.Cells.Item(51,5).Value = '  +4.66%  '

# Restore default cell style (values remain text) to match original formatting
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(51,4).Style = "Normal"
